{"js": "// 1) Update the \"last updated\" date from 2022-06-28 to 2022-11-04.\nconst dateHits = context.document.body.search(\"2022-06-28\", { matchCase: true });\ndateHits.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < dateHits.items.length; i++) {\n  dateHits.items[i].insertText(\"2022-11-04\", \"Replace\");\n}\nawait context.sync();\n\n// 2) Remove the six \"Source Code\" style paragraphs (package-build warnings)\n//    that sit right after the \"Connors et al. (2020).\" sentence.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst items = paragraphs.items;\nitems.forEach((p) => p.load(\"text,style\"));\nawait context.sync();\n\nconst toDelete = items.filter(\n  (p) => p.style === \"Source Code\" && p.text.indexOf(\"was built under R version\") !== -1\n);\ntoDelete.forEach((p) => p.delete());\nawait context.sync();\n\n// 3) Fix the typographic multiplication sign in the tibble dimensions line.\nconst tibbleHits = context.document.body.search(\"## # A tibble: 3 x 6\", { matchCase: true });\ntibbleHits.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < tibbleHits.items.length; i++) {\n  tibbleHits.items[i].insertText(\"## # A tibble: 3 \u00d7 6\", \"Replace\");\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Update the \"last updated\" date from 2022-06-28 to 2022-11-04.\n$find1 = $d.Content\n$find1.Find.Execute(\n    \"2022-06-28\", $false, $false, $false, $false, $false,\n    $true, 1, $false, \"2022-11-04\", 2\n) | Out-Null\n\n# 2) Remove the six \"Source Code\" style paragraphs (package-build warnings)\n#    that sit right after the \"Connors et al. (2020).\" sentence. Walk\n#    backwards so deleting doesn't shift the indices of paragraphs still to\n#    be inspected.\nfor ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {\n    $p = $d.Paragraphs.Item($i)\n    $style = $p.Range.Style.NameLocal\n    $text = $p.Range.Text\n    if ($style -eq \"Source Code\" -and $text -like \"*was built under R version*\") {\n        $p.Range.Delete()\n    }\n}\n\n# 3) Fix the typographic multiplication sign in the tibble dimensions line.\n$find2 = $d.Content\n$find2.Find.Execute(\n    \"## # A tibble: 3 x 6\", $false, $false, $false, $false, $false,\n    $true, 1, $false, \"## # A tibble: 3 \u00d7 6\", 2\n) | Out-Null\n"}
